# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain numeric-looking strings (e.g. "398.36").
# Mark those specific cells as Text first so Excel keeps them as literal strings
# (matching the source data) instead of auto-converting them to numbers.
# (looping cell-by-cell rather than using a comma union Range, which only
# applies NumberFormat to the first area)
$textPriceCells = @("D5","D6","D7","D9","D10","D11","D14","D15","D17","D18","D22","D23","D24","D25","D26","D27","D28","D29","D31","D33","D34","D36","D37","D38","D40","D41","D45","D46","D47","D49")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Row-by-row updates -------------------------------------------------
# Row 2
$ws.Range("D2").Value = '56.398.04'
$ws.Range("E2").Value = '  +10.22%  '

# Row 3
$ws.Range("D3").Value = '3.245.62'
$ws.Range("E3").Value = '  +6.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '398.36'
$ws.Range("E5").Value = '  +3.16%  '

# Row 6
$ws.Range("D6").Value = '111.33'
$ws.Range("E6").Value = '  +8.58%  '

# Row 7
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  +3.93%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  +7.20%  '

# Row 10
$ws.Range("D10").Value = '39.33'
$ws.Range("E10").Value = '  +7.26%  '

# Row 11
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  +10.09%  '

# Row 12
$ws.Range("E12").Value = '  +2.33%  '

# Row 13
$ws.Range("D13").Value = '3.739.77'
$ws.Range("E13").Value = '  +4.97%  '

# Row 14
$ws.Range("D14").Value = '19.19'
$ws.Range("E14").Value = '  +4.75%  '

# Row 15
$ws.Range("D15").Value = '8.12'
$ws.Range("E15").Value = '  +5.65%  '

# Row 16
$ws.Range("D16").Value = '3.229.94'
$ws.Range("E16").Value = '  +5.85%  '

# Row 17
$ws.Range("D17").Value = '1.05'
$ws.Range("E17").Value = '  +6.09%  '

# Row 18
$ws.Range("D18").Value = '11.07'
$ws.Range("E18").Value = '  +4.10%  '

# Row 19
$ws.Range("D19").Value = '56.250.77'
$ws.Range("E19").Value = '  +9.86%  '

# Row 20
$ws.Range("E20").Value = '  +4.44%  '

# Row 21
$ws.Range("E21").Value = '  +8.58%  '

# Row 22
$ws.Range("D22").Value = '13.12'
$ws.Range("E22").Value = '  +6.68%  '

# Row 23
$ws.Range("D23").Value = '298.14'
$ws.Range("E23").Value = '  +12.82%  '

# Row 24
$ws.Range("D24").Value = '75.84'
$ws.Range("E24").Value = '  +8.76%  '

# Row 25
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +3.33%  '

# Row 26
$ws.Range("D26").Value = '8.21'
$ws.Range("E26").Value = '  +3.85%  '

# Row 27
$ws.Range("D27").Value = '28.21'
$ws.Range("E27").Value = '  +4.74%  '

# Row 28
$ws.Range("D28").Value = '7.38'
$ws.Range("E28").Value = '  +2.59%  '

# Row 29
$ws.Range("D29").Value = '0.171'
$ws.Range("E29").Value = '  +4.46%  '

# Row 30
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("D31").Value = '0.111'
$ws.Range("E31").Value = '  +5.31%  '

# Row 32
$ws.Range("E32").Value = '  +6.94%  '

# Row 33
$ws.Range("D33").Value = '37.06'
$ws.Range("E33").Value = '  +4.54%  '

# Row 34
$ws.Range("D34").Value = '0.0492'
$ws.Range("E34").Value = '  +4.43%  '

# Row 35
$ws.Range("E35").Value = '  +3.42%  '

# Row 36
$ws.Range("D36").Value = '51.44'
$ws.Range("E36").Value = '  +3.33%  '

# Row 37
$ws.Range("D37").Value = '3.56'

# Row 38
$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  +25.54%  '

# Row 39
$ws.Range("E39").Value = '  -0.13%  '

# Row 40
$ws.Range("D40").Value = '134.91'
$ws.Range("E40").Value = '  +4.55%  '

# Row 41
$ws.Range("D41").Value = '17.58'
$ws.Range("E41").Value = '  +6.35%  '

# Row 42
$ws.Range("E42").Value = '  +5.43%  '

# Row 43
$ws.Range("E43").Value = '  +4.62%  '

# Row 44
$ws.Range("E44").Value = '  +4.21%  '

# Row 45
$ws.Range("D45").Value = '0.284'
$ws.Range("E45").Value = '  -1.23%  '

# Row 46
$ws.Range("D46").Value = '22.42'
$ws.Range("E46").Value = '  +3.48%  '

# Row 47
$ws.Range("D47").Value = '2.20'
$ws.Range("E47").Value = '  +57.39%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.139.57'
$ws.Range("E48").Value = '  +3.93%  '

# Row 49
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '2.09'
$ws.Range("E49").Value = '  +0.35%  '

# Row 50
$ws.Range("E50").Value = '  -1.71%  '

# Row 51
$ws.Range("E51").Value = '  +12.15%  '
